$wb = $excel.ActiveWorkbook

$wsDashboard = $wb.Worksheets.Item("DashboardPage")
$wsPim       = $wb.Worksheets.Item("PIMPage")

# --- Update existing Employee Id locator on PIMPage (row 5) ---
$wsPim.Range("C5").Value = "//label[text()='Employee Id']/parent::div/following-sibling::div/input"

# --- Add new "Search Employee" related locators on PIMPage ---
$wsPim.Range("A12").Value = "searchEmployeeNameField"
$wsPim.Range("B12").Value = "XPATH"

$wsPim.Range("A13").Value = "searchEmployeeButtonLocator"
$wsPim.Range("B13").Value = "XPATH"

$wsPim.Range("A16").Value = "searchedEmployeeIdRecordTable"
$wsPim.Range("B16").Value = "XPATH"

# --- Add the profile-name locator on DashboardPage (row 3) ---
$wsDashboard.Range("A3").Value = "profileNameLocator"
$wsDashboard.Range("B3").Value = "XPATH"

# --- Fill in the XPATH values for the new Search Employee rows ---
$wsPim.Range("C12").Value = "//label[text()='Employee Name']/parent::div/following-sibling::div/div/div/input"
$wsPim.Range("C13").Value = "//button[text()=' Search ']"

$wsPim.Range("A14").Value = "searchedEmployeeFirstNameRecordTable"
$wsPim.Range("B14").Value = "XPATH"

$wsPim.Range("A15").Value = "searchedEmployeeLastNameRecordTable"
$wsPim.Range("B15").Value = "XPATH"

$wsPim.Range("C14").Value = "(//div[@class='oxd-table-cell oxd-padding-cell']/div)[3]"
$wsPim.Range("C15").Value = "(//div[@class='oxd-table-cell oxd-padding-cell']/div)[4]"
$wsPim.Range("C16").Value = "(//div[@class='oxd-table-cell oxd-padding-cell']/div)[2]"

$wsDashboard.Range("C3").Value = "//p[@class='oxd-userdropdown-name']"

# --- Column width tweak on PIMPage column C (best effort) ---
$wsPim.Columns.Item(3).ColumnWidth = 67.5

# --- Selections / active sheet to match final saved view state ---
$wsPim.Range("C15").Select()
$wsDashboard.Activate()
$wsDashboard.Range("C3").Select()
